$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "42.811.73"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.291.24"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.99"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.01"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.87"
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.76"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "2.648.11"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "2.292.20"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "42.752.06"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.53"
$ws.Range("E19").Value = "  -4.68%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.78"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.92"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.15"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.18"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.85"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.13"
$ws.Range("E36").Value = "  -6.82%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0685"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.109"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "2.012.68"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.08"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.09"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("E48").Value = "  -2.56%  "
$ws.Range("D49").Value = "2.515.83"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.99"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.77"
$ws.Range("E51").Value = "  -8.06%  "
